$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set value and custom number format (Hryvnia currency-like format) on B4,
# matching the new numFmt/cellXfs entries added to styles.xml
$cell = $ws.Range("B4")
$cell.Value = 7890.1234560000003
$cell.NumberFormat = "\₴#,###"

# Update the active selection to B4 (was A7)
$ws.Range("B4").Select()
